# observaciones.xlsx update
# - appends a new observation (row 7) to the "15-08-2025" sheet (2nd tab)
# - adds a brand new day sheet "11-08-2025" with its header + one observation row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Append new observation row to the existing "15-08-2025" sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Copy the formatting of the last existing data row (row 6) onto the new row (7)
# so the new cells pick up the same style index (header=1 / data=2) used
# throughout the workbook instead of the default style.
$ws2.Range("A6:F6").Copy()
$ws2.Range("A7:F7").PasteSpecial(-4122)

$ws2.Cells.Item(7, 1).Value = 12
$ws2.Cells.Item(7, 2).Value = "23:04"
$ws2.Cells.Item(7, 3).Value = "LÍNEA 2"
$ws2.Cells.Item(7, 4).Value = "MÁQUINA 29 T8"
$ws2.Cells.Item(7, 5).Value = "[General] a"
$ws2.Cells.Item(7, 6).Value = "admin"

# ---------------------------------------------------------------------------
# 2) Add a new sheet "11-08-2025" at the end of the workbook with one
#    observation (ID 13)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "11-08-2025"

# Re-use the header/data formatting (style indexes 1 and 2) already defined
# in the workbook by pasting formats from the "15-08-2025" sheet.
$wb.Worksheets.Item(2).Range("A1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)

$wb.Worksheets.Item(2).Range("A6").Copy()
$newSheet.Range("A2:F2").PasteSpecial(-4122)

$newSheet.Cells.Item(1, 1).Value = "ID"
$newSheet.Cells.Item(1, 2).Value = "Hora"
$newSheet.Cells.Item(1, 3).Value = "Línea"
$newSheet.Cells.Item(1, 4).Value = "Máquina"
$newSheet.Cells.Item(1, 5).Value = "Observación"
$newSheet.Cells.Item(1, 6).Value = "Usuario"

$newSheet.Cells.Item(2, 1).Value = 13
$newSheet.Cells.Item(2, 2).Value = "23:21"
$newSheet.Cells.Item(2, 3).Value = "LÍNEA 2"
$newSheet.Cells.Item(2, 4).Value = "MÁQUINA 28 T12"
$newSheet.Cells.Item(2, 5).Value = "[General] aaa"
$newSheet.Cells.Item(2, 6).Value = "admin"

# Column widths for the new sheet (8,12,15,20,50,15 characters). The
# ColumnWidth setter in this runtime stores a small fixed pixel padding on
# top of the requested character width, so the value is pre-compensated
# here to land exactly on the intended integer width.
$wOffset = 0.83333333333333
$newSheet.Columns.Item(1).ColumnWidth = 8 - $wOffset
$newSheet.Columns.Item(2).ColumnWidth = 12 - $wOffset
$newSheet.Columns.Item(3).ColumnWidth = 15 - $wOffset
$newSheet.Columns.Item(4).ColumnWidth = 20 - $wOffset
$newSheet.Columns.Item(5).ColumnWidth = 50 - $wOffset
$newSheet.Columns.Item(6).ColumnWidth = 15 - $wOffset

$newSheet.Range("A1").Select()
